$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.378.98"
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").Value = "1.596.10"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "1.824.98"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").Value = "1.596.99"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.534"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "28.395.24"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.22%  "
$ws.Range("D19").Value = "0.0₃0711"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.107"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "1.398.97"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.64%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.73%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.540"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "1.734.89"
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0526"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
